# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 12:43"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 4433532
$ws.Range("C4").Value = 122
$ws.Range("D4").Value = 2137187
$ws.Range("E4").Value = 2145895
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 150450

# Row 14: Iran
$ws.Range("B14").Value = 296273
$ws.Range("C14").Value = 2667
$ws.Range("D14").Value = 257019
$ws.Range("E14").Value = 23107
$ws.Range("G14").Value = 235
$ws.Range("H14").Value = 16147

# Row 21: Alemania
$ws.Range("B21").Value = 207416
$ws.Range("C21").Value = 37
$ws.Range("E21").Value = 6810
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 9206

# Row 43: Emiratos Arabes Unidos
$ws.Range("B43").Value = 59546
$ws.Range("C43").Value = 369
$ws.Range("D43").Value = 52905
$ws.Range("E43").Value = 6294
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 347

# Row 47: Rumania
$ws.Range("B47").Value = 47053
$ws.Range("C47").Value = 1151
$ws.Range("D47").Value = 26128
$ws.Range("E47").Value = 18686
$ws.Range("G47").Value = 33
$ws.Range("H47").Value = 2239

# Row 74: was Australia, now El Salvador (El Salvador overtook Australia in ranking)
$ws.Range("A74").Value = "El Salvador"
$ws.Range("B74").Value = 15446
$ws.Range("C74").Value = 411
$ws.Range("D74").Value = 7903
$ws.Range("E74").Value = 7126
$ws.Range("G74").Value = 9
$ws.Range("H74").Value = 417

# Row 75: was El Salvador, now Australia
$ws.Range("A75").Value = "Australia"
$ws.Range("B75").Value = 15302
$ws.Range("C75").Value = 367
$ws.Range("D75").Value = 9311
$ws.Range("E75").Value = 5824
$ws.Range("G75").Value = 6
$ws.Range("H75").Value = 167

# Row 84: Senegal
$ws.Range("B84").Value = 9805
$ws.Range("C84").Value = 41
$ws.Range("D84").Value = 6591
$ws.Range("E84").Value = 3016
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 198
